$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 775
$ws.Range("F5").Value = 2265
$ws.Range("F6").Value = 1347
$ws.Range("F7").Value = 102
$ws.Range("F9").Value = 129
$ws.Range("F11").Value = 2955
$ws.Range("F19").Value = 1018
$ws.Range("F20").Value = 1018
$ws.Range("F25").Value = 187
$ws.Range("F31").Value = 1014
$ws.Range("F32").Value = 4990
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F34").Value = 319
$ws.Range("F35").Value = 319
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 775
$ws.Range("F10").Value = 2265
$ws.Range("F11").Value = 1347
$ws.Range("F12").Value = 102
$ws.Range("F14").Value = 129
$ws.Range("F18").Value = 2955
$ws.Range("F27").Value = 1018
$ws.Range("F28").Value = 1018
$ws.Range("F34").Value = 187
$ws.Range("F45").Value = 1014
$ws.Range("F46").Value = 4990
$ws.Range("F49").Value = 319
